$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 2-39: update Price (D) and Volume(1h) (E) columns.
# Numeric-looking Price values are assigned with a leading apostrophe so
# Excel keeps them as text (matching the original inlineStr cell type)
# instead of auto-converting them to numbers.
$ws.Range("D2").Value = "26.096.79"
$ws.Range("E2").Value = "  -0.24%  "
$ws.Range("D3").Value = "1.647.88"
$ws.Range("E3").Value = "  -1.30%  "
$ws.Range("D4").Value = "'1.001"
$ws.Range("E4").Value = "  -0.12%  "
$ws.Range("D5").Value = "'216.85"
$ws.Range("E5").Value = "  +3.32%  "
$ws.Range("D6").Value = "'0.5216"
$ws.Range("E6").Value = "  -0.27%  "
$ws.Range("D7").Value = "'1.002"
$ws.Range("E7").Value = "  -0.05%  "
$ws.Range("D8").Value = "'0.2610"
$ws.Range("E8").Value = "  -0.52%  "
$ws.Range("D9").Value = "'0.06459"
$ws.Range("E9").Value = "  +2.04%  "
$ws.Range("D10").Value = "'20.87"
$ws.Range("E10").Value = "  -1.64%  "
$ws.Range("D11").Value = "'0.07702"
$ws.Range("E11").Value = "  +2.08%  "
$ws.Range("D12").Value = "1.666.17"
$ws.Range("E12").Value = "  -0.32%  "
$ws.Range("D13").Value = "'4.425"
$ws.Range("E13").Value = "  -0.46%  "
$ws.Range("D14").Value = "1.869.12"
$ws.Range("E14").Value = "  -1.62%  "
$ws.Range("D15").Value = "'0.5581"
$ws.Range("E15").Value = "  +1.35%  "
$ws.Range("D16").Value = "0.0₅8296"
$ws.Range("E16").Value = "  +3.65%  "
$ws.Range("D17").Value = "'65.22"
$ws.Range("E17").Value = "  -2.07%  "
$ws.Range("D18").Value = "26.100.73"
$ws.Range("E19").Value = "  +0.01%  "
$ws.Range("D20").Value = "'4.743"
$ws.Range("E20").Value = "  -0.45%  "
$ws.Range("D21").Value = "'188.61"
$ws.Range("D22").Value = "'10.24"
$ws.Range("E22").Value = "  -0.99%  "
$ws.Range("D23").Value = "'6.226"
$ws.Range("E23").Value = "  +0.33%  "
$ws.Range("D24").Value = "'1.002"
$ws.Range("E24").Value = "  -0.12%  "
$ws.Range("D25").Value = "'146.36"
$ws.Range("E25").Value = "  -2.23%  "
$ws.Range("D26").Value = "'7.452"
$ws.Range("E26").Value = "  -0.77%  "
$ws.Range("D27").Value = "'0.1213"
$ws.Range("E27").Value = "  -2.87%  "
$ws.Range("D28").Value = "'15.85"
$ws.Range("E28").Value = "  -0.16%  "
$ws.Range("D29").Value = "'1.401"
$ws.Range("E29").Value = "  +3.73%  "
$ws.Range("E30").Value = "  -7.21%  "
$ws.Range("D31").Value = "'1.269"
$ws.Range("E31").Value = "  -0.49%  "
$ws.Range("D32").Value = "'3.420"
$ws.Range("E32").Value = "  -2.68%  "
$ws.Range("E33").Value = "  -0.24%  "
$ws.Range("D34").Value = "'1.656"
$ws.Range("E34").Value = "  +0.46%  "
$ws.Range("D35").Value = "'0.9905"
$ws.Range("E35").Value = "  -1.46%  "
$ws.Range("E36").Value = "  -0.58%  "
$ws.Range("E37").Value = "  +0.14%  "
$ws.Range("D38").Value = "'0.5663"
$ws.Range("E38").Value = "  -6.15%  "
$ws.Range("D39").Value = "'0.01619"
$ws.Range("E39").Value = "  +0.30%  "

# Rows 40/41: the two coins swap places (TrustWalletToken <-> FraxShare),
# and both get new Price/Volume values.
$ws.Range("B40").Value = "FraxShare"
$ws.Range("C40").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D40").Value = "'5.848"
$ws.Range("E40").Value = "  -4.80%  "

$ws.Range("B41").Value = "TrustWalletToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D41").Value = "'0.8570"
$ws.Range("E41").Value = "  -1.04%  "

# Rows 42-51: update remaining Price (D) and Volume(1h) (E) columns.
$ws.Range("E42").Value = "  -0.19%  "
$ws.Range("D43").Value = "1.032.66"
$ws.Range("E43").Value = "  -7.06%  "
$ws.Range("D44").Value = "'100.31"
$ws.Range("E44").Value = "  -0.03%  "
$ws.Range("D45").Value = "1.797.44"
$ws.Range("E45").Value = "  -1.40%  "
$ws.Range("D46").Value = "0.0₈108"
$ws.Range("E46").Value = "  +0.15%  "
$ws.Range("D47").Value = "'56.03"
$ws.Range("E47").Value = "  +0.89%  "
$ws.Range("D48").Value = "'1.004"
$ws.Range("E48").Value = "  +0.11%  "
$ws.Range("D49").Value = "'8.080"
$ws.Range("E49").Value = "  +0.34%  "
$ws.Range("E50").Value = "  -0.92%  "
$ws.Range("D51").Value = "'0.4221"
$ws.Range("E51").Value = "  -0.42%  "
